# Generate Report for Handback
# Updates the localization-status workbook to reflect that both language
# reports (zh-cn, de-de) have been handed back: status text, the "Latest
# Target File" / "Latest Handback File" columns are populated (with a
# hyperlink on the target-file cell, mirroring the source-file hyperlink),
# the handback datetime is recorded, and the now-wider columns are resized.

$wb = $excel.ActiveWorkbook

$urlDoc1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f488bdce8bcd488f326c5196d7345f94b018958e/e2e/84243bd6-a7a3-4bf9-ac01-3b3c2383be3f.md"
$urlDoc2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f488bdce8bcd488f326c5196d7345f94b018958e/e2e/d859c7c7-4e5b-4868-8128-c84ee998d033.md"

$doc1Name = "84243bd6-a7a3-4bf9-ac01-3b3c2383be3f.md"
$doc2Name = "d859c7c7-4e5b-4868-8128-c84ee998d033.md"

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (E/F) now read "Handed back", and those
# two columns are widened to fit the longer text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn report
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C3").Value = $statusHandedBack

$zhcn.Range("I2").Value = $doc1Name
$zhcn.Range("J2").Value = "84243bd6-a7a3-4bf9-ac01-3b3c2383be3f.6bf8cfdd963cf7ab558a8d355ffe6ddc51ed15e6.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-25 20:27:49"

$zhcn.Range("I3").Value = $doc2Name
$zhcn.Range("J3").Value = "d859c7c7-4e5b-4868-8128-c84ee998d033.1dbd9ed04fb5a636f04db86c3416bb1d08f6d7f8.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-25 20:27:49"

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# Add hyperlinks on the new "Latest Target File" cells (I2/I3) linking to
# the same source document as the existing A2/A3 cells (left untouched).
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $urlDoc1, "", "", $doc1Name)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urlDoc2, "", "", $doc2Name)

# ---------------------------------------------------------------------
# de-de report
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C3").Value = $statusHandedBack

$dede.Range("I2").Value = $doc1Name
$dede.Range("J2").Value = "84243bd6-a7a3-4bf9-ac01-3b3c2383be3f.6bf8cfdd963cf7ab558a8d355ffe6ddc51ed15e6.de-de.xlf"
$dede.Range("K2").Value = "2016-08-25 20:27:56"

$dede.Range("I3").Value = $doc2Name
$dede.Range("J3").Value = "d859c7c7-4e5b-4868-8128-c84ee998d033.1dbd9ed04fb5a636f04db86c3416bb1d08f6d7f8.de-de.xlf"
$dede.Range("K3").Value = "2016-08-25 20:27:56"

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $urlDoc1, "", "", $doc1Name)
$dede.Hyperlinks.Add($dede.Range("I2"), $urlDoc1, "", "", $doc1Name)
$dede.Hyperlinks.Add($dede.Range("A3"), $urlDoc2, "", "", $doc2Name)
$dede.Hyperlinks.Add($dede.Range("I3"), $urlDoc2, "", "", $doc2Name)
